$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bảng TinRaoVat: thêm cột GhiChu -> đánh dấu "X" hoàn thành
# cho các mục tương ứng tại cột E (Hoàn thành), dòng 2 đến 6
$ws.Range("E2:E6").Value = "X"
